$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $escaped = $value -replace '"', '""'
    $ws.Range($cellRef).Formula = '="' + $escaped + '"'
    $ws.Range($cellRef).Copy()
    $ws.Range($cellRef).PasteSpecial(-4163)
}

Set-TextValue 'D2' '28.921.85'
Set-TextValue 'E2' '  -2.41%  '
Set-TextValue 'D3' '1.902.20'
Set-TextValue 'E3' '  -4.58%  '
Set-TextValue 'D4' '1.005'
Set-TextValue 'E4' '  +0.18%  '
Set-TextValue 'D5' '324.08'
Set-TextValue 'E5' '  -1.89%  '
Set-TextValue 'D6' '1.003'
Set-TextValue 'E6' '  +0.05%  '
Set-TextValue 'D7' '0.4593'
Set-TextValue 'E7' '  -2.14%  '
Set-TextValue 'D8' '0.3811'
Set-TextValue 'E8' '  -3.73%  '
Set-TextValue 'B9' 'OKB'
Set-TextValue 'C9' 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextValue 'D9' '45.55'
Set-TextValue 'E9' '  -2.32%  '
Set-TextValue 'B10' 'Dogecoin'
Set-TextValue 'C10' 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
Set-TextValue 'D10' '0.07713'
Set-TextValue 'E10' '  -4.66%  '
Set-TextValue 'B11' 'Polygon'
Set-TextValue 'C11' 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
Set-TextValue 'D11' '0.9808'
Set-TextValue 'E11' '  -2.24%  '
Set-TextValue 'B12' 'Solana'
Set-TextValue 'C12' 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
Set-TextValue 'D12' '22.05'
Set-TextValue 'E12' '  -3.97%  '
Set-TextValue 'B13' 'WrappedEther'
Set-TextValue 'C13' 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextValue 'D13' '1.980.55'
Set-TextValue 'E13' '  -0.79%  '
Set-TextValue 'B14' 'Chainlink'
Set-TextValue 'C14' 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextValue 'D14' '6.963'
Set-TextValue 'E14' '  -4.11%  '
Set-TextValue 'B15' 'Polkadot'
Set-TextValue 'C15' 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextValue 'D15' '5.669'
Set-TextValue 'E15' '  -3.59%  '
Set-TextValue 'B16' 'TRON'
Set-TextValue 'C16' 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
Set-TextValue 'D16' '0.07051'
Set-TextValue 'E16' '  -1.37%  '
Set-TextValue 'B17' 'BinanceUSD'
Set-TextValue 'C17' 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
Set-TextValue 'D17' '1.005'
Set-TextValue 'E17' '  -0.02%  '
Set-TextValue 'B18' 'Litecoin'
Set-TextValue 'C18' 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
Set-TextValue 'D18' '84.01'
Set-TextValue 'E18' '  -5.58%  '
Set-TextValue 'B19' 'ShibaInu'
Set-TextValue 'C19' 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
Set-TextValue 'D19' '0.000009535'
Set-TextValue 'E19' '  -5.18%  '
Set-TextValue 'B20' 'Avalanche'
Set-TextValue 'C20' 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
Set-TextValue 'D20' '16.71'
Set-TextValue 'E20' '  -4.28%  '
Set-TextValue 'B21' 'Dai'
Set-TextValue 'C21' 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextValue 'D21' '1.003'
Set-TextValue 'E21' '  +0.09%  '
Set-TextValue 'B22' 'WrappedBTC'
Set-TextValue 'C22' 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
Set-TextValue 'D22' '28.905.21'
Set-TextValue 'E22' '  -2.49%  '
Set-TextValue 'B23' 'Uniswap'
Set-TextValue 'C23' 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
Set-TextValue 'D23' '5.324'
Set-TextValue 'E23' '  -4.36%  '
Set-TextValue 'B24' 'Cosmos'
Set-TextValue 'C24' 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextValue 'D24' '10.89'
Set-TextValue 'E24' '  -3.50%  '
Set-TextValue 'B25' 'WrappedliquidstakedEther2.0'
Set-TextValue 'C25' 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
Set-TextValue 'D25' '2.181.43'
Set-TextValue 'E25' '  -2.48%  '
Set-TextValue 'B26' 'Toncoin'
Set-TextValue 'C26' 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextValue 'D26' '2.095'
Set-TextValue 'E26' '  -1.14%  '
Set-TextValue 'B27' 'Monero'
Set-TextValue 'C27' 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue 'D27' '157.47'
Set-TextValue 'E27' '  -0.31%  '
Set-TextValue 'B28' 'EthereumClassic'
Set-TextValue 'C28' 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextValue 'D28' '19.15'
Set-TextValue 'E28' '  -2.97%  '
Set-TextValue 'B29' 'InternetComputer(DFINITY)'
Set-TextValue 'C29' 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextValue 'D29' '5.589'
Set-TextValue 'E29' '  -7.18%  '
Set-TextValue 'B30' 'BitcoinCash'
Set-TextValue 'C30' 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
Set-TextValue 'D30' '117.61'
Set-TextValue 'E30' '  -2.33%  '
Set-TextValue 'B31' 'LidoDAOToken'
Set-TextValue 'C31' 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
Set-TextValue 'D31' '1.853'
Set-TextValue 'E31' '  -4.83%  '
Set-TextValue 'B32' 'Stellar'
Set-TextValue 'C32' 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue 'D32' '0.09283'
Set-TextValue 'E32' '  -1.96%  '
Set-TextValue 'B33' 'ImmutableX'
Set-TextValue 'C33' 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue 'D33' '0.8622'
Set-TextValue 'E33' '  -6.10%  '
Set-TextValue 'B34' 'Filecoin'
Set-TextValue 'C34' 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue 'D34' '5.093'
Set-TextValue 'E34' '  -3.68%  '
Set-TextValue 'B35' 'ARBITRUM'
Set-TextValue 'C35' 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextValue 'D35' '1.253'
Set-TextValue 'E35' '  -7.63%  '
Set-TextValue 'B36' 'HuobiToken'
Set-TextValue 'C36' 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
Set-TextValue 'D36' '3.014'
Set-TextValue 'E36' '  -5.38%  '
Set-TextValue 'B37' 'Hedera'
Set-TextValue 'C37' 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue 'D37' '0.05695'
Set-TextValue 'E37' '  -2.83%  '
Set-TextValue 'B38' 'TrustWalletToken'
Set-TextValue 'C38' 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextValue 'D38' '1.147'
Set-TextValue 'E38' '  -2.39%  '
Set-TextValue 'B39' 'Frax'
Set-TextValue 'C39' 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
Set-TextValue 'D39' '1.003'
Set-TextValue 'E39' '  -0.15%  '
Set-TextValue 'B40' 'VeChain'
Set-TextValue 'C40' 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue 'D40' '0.02038'
Set-TextValue 'E40' '  -4.40%  '
Set-TextValue 'B41' 'FraxShare'
Set-TextValue 'C41' 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue 'D41' '7.484'
Set-TextValue 'E41' '  -5.75%  '
Set-TextValue 'B42' 'TheSandbox'
Set-TextValue 'C42' 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
Set-TextValue 'D42' '0.5511'
Set-TextValue 'E42' '  -4.91%  '
Set-TextValue 'B43' 'Algorand'
Set-TextValue 'C43' 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextValue 'D43' '0.1753'
Set-TextValue 'E43' '  -4.33%  '
Set-TextValue 'B44' 'Aptos'
Set-TextValue 'C44' 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextValue 'D44' '9.305'
Set-TextValue 'E44' '  -6.23%  '
Set-TextValue 'B45' 'MXToken'
Set-TextValue 'C45' 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextValue 'D45' '2.731'
Set-TextValue 'E45' '  -1.17%  '
Set-TextValue 'B46' 'Decentraland'
Set-TextValue 'C46' 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
Set-TextValue 'D46' '0.5194'
Set-TextValue 'E46' '  -3.80%  '
Set-TextValue 'B47' 'EnergySwap'
Set-TextValue 'C47' 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue 'D47' '11.35'
Set-TextValue 'E47' '  -6.42%  '
Set-TextValue 'B48' 'RenderToken'
Set-TextValue 'C48' 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue 'D48' '2.081'
Set-TextValue 'E48' '  -6.05%  '
Set-TextValue 'B49' 'Cronos'
Set-TextValue 'C49' 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextValue 'D49' '0.06820'
Set-TextValue 'E49' '  -2.09%  '
Set-TextValue 'B50' 'Quant'
Set-TextValue 'C50' 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
Set-TextValue 'D50' '111.20'
Set-TextValue 'E50' '  -2.70%  '
Set-TextValue 'B51' 'NEARProtocol'
Set-TextValue 'C51' 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextValue 'D51' '1.773'
Set-TextValue 'E51' '  -5.54%  '
